# Separate suites for qa and stging and message update
#
# Updates the product-code placeholders for the three rows on the
# "Input" sheet that carry a distinct per-row style (B2, B5, B8).
# Each cell keeps its existing formatting; only the displayed
# product code text changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "prodcxsf"
$ws.Range("B5").Value = "prodnfTi"
$ws.Range("B8").Value = "prodreFi"
